$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Odd_H_FT","Odd_D_FT","Odd_A_FT","Odd_Over15_FT","Odd_Under15_FT",
    "Odd_Over25_FT","Odd_Under25_FT","Odd_Over05_HT","Odd_Under05_HT",
    "Odd_BTTS_Yes","Odd_BTTS_No","Odd_CS_1-0","Odd_CS_2-0","Odd_CS_2-1",
    "Odd_CS_3-0","Odd_CS_3-1","Odd_CS_3-2","Odd_CS_0-0","Odd_CS_1-1",
    "Odd_CS_2-2","Odd_CS_3-3","Odd_CS_4-4","Odd_CS_0-1","Odd_CS_0-2",
    "Odd_CS_1-2","Odd_CS_0-3","Odd_CS_1-3","Odd_CS_2-3"
)

$row2values = @(2.67,3.2,2.47,1.4,2.55,2.15,1.55,1.47,2.32,1.91,1.7,7.2,12,10.5,30,26,40,7.9,6.2,17,100,900,6.9,11,10,26,23,40)

$startCol = 7  # column G

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $startCol + $i

    # Row 1 - header labels
    $cell1 = $ws.Cells.Item(1, $col)
    $cell1.Value = $headers[$i]

    # Row 2 - odds values
    $cell2 = $ws.Cells.Item(2, $col)
    $cell2.Value = $row2values[$i]

    # Row 3 - blank placeholder cell (forces the cell to exist without a value)
    $cell3 = $ws.Cells.Item(3, $col)
    $cell3.Borders.LineStyle = -4142
}

# Apply the existing header style (from A1) to the newly added header cells
$srcHeader = $ws.Range("A1")
$srcHeader.Copy()
$destHeaders = $ws.Range("G1:AH1")
$destHeaders.PasteSpecial(-4122)
$excel.CutCopyMode = 0
